# ДК.xlsx — schedule update: refresh stage names/owners/dates, drop 2 filler rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two blank filler rows (A11, A12) — shifts rows 13:20 up to 11:18.
$ws.Rows("11:12").Delete()

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = 45386
$ws.Cells.Item(2,2).Value = "Создаем веб-страницу с регистрацией"
$ws.Cells.Item(2,4).Value = "Никита"

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = 45392
$ws.Cells.Item(3,2).Value = "Страница с рулеткой"
$ws.Cells.Item(3,4).Value = "Никита"

# --- Row 4 ---
$ws.Cells.Item(4,1).Value = 45395
$ws.Cells.Item(4,2).Value = "Страница с интуицией "
$ws.Cells.Item(4,4).Value = "Тимур"

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = 45398
$ws.Cells.Item(5,2).Value = "Страница вывода"
$ws.Cells.Item(5,4).Value = "Аня"

# --- Row 6 ---
$ws.Cells.Item(6,1).Value = 45400
$ws.Cells.Item(6,2).Value = "Страница депа"
$ws.Cells.Item(6,4).Value = "Аня"

# --- Row 7 ---
$ws.Cells.Item(7,1).Value = 45402
$ws.Cells.Item(7,2).Value = "Обработка денег"
$ws.Cells.Item(7,4).Value = "Никита"

# --- Row 8 ---
$ws.Cells.Item(8,1).Value = 45403
$ws.Cells.Item(8,2).Value = "Бот"
$ws.Cells.Item(8,4).Value = "Тимур"

# --- Row 9 ---
$ws.Cells.Item(9,1).Value = 45404
$ws.Cells.Item(9,2).Value = "Все готово"
$ws.Cells.Item(9,4).Value = "Все участники"

# --- Row 10 ---
$ws.Cells.Item(10,1).Value = 45407
$ws.Cells.Item(10,2).Value = "Защита"
$ws.Cells.Item(10,4).Value = "Все участники"

# Match the final on-screen selection left by the edit.
$ws.Range("A9:D10").Select()
